# Update column C (Förändrad) for rows 2-43 from 45758 to 45759 (increment by 1 day)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45758) {
        $cell.Value = 45759
    }
}
